# Se incorporan los dos test restantes y modificaciones varias
# Fill in the two remaining test cases (rows 3 and 4) of the
# "Detalle_tests" sheet, resize column C, adjust row heights, and
# move the active selection to H3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "Chequear si se puede crear un post correctamente" ---
$ws.Range("B3").Value = "Chequear si se puede crear un post correctamente"
$ws.Range("C3").Value = "titulo=random`nequipo=random`nposteo=`"Esta es una prueba 1234567890 /&%#$`"`nautor_prueba=User.objects.create(id=1)"
$ws.Range("D3").Value = "El posteo se debe crear correctamente "
$ws.Range("E3").Value = "El posteo se creó en el blog correctamente"
$ws.Range("F3").Value = "Pass"

# --- Row 4: "Corroborar si toma como válido el registro de usuario" ---
$ws.Range("B4").Value = "Corroborar si toma como válido el registro de usuario"
# Leading apostrophe forces a text "quote-prefix" cell (the value itself keeps
# its own leading space, just like it was typed into Excel as
# `' 'last_name': ...`)
$ws.Range("C4").Value = "' 'last_name': 'Ronaldo',`n'first_name': 'Cristiano',`n'username': 'cronaldo7',`n'email':'cronaldo7@example.com',`n'password1': 'Portugal123',`n'password2': 'Portugal123'"
$ws.Range("D4").Value = "El registro debe ser válido"
$ws.Range("E4").Value = "El registro fue válido"
$ws.Range("F4").Value = "Pass"

# Column width + row heights
$ws.Columns("C").ColumnWidth = 41.5
$ws.Rows(3).RowHeight = 61.5
$ws.Rows(4).RowHeight = 94.5

# Move selection
$ws.Range("H3").Select()
